$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-10: full ligand-receptor data for the ECs / FAPs / Neutro 3x3 cluster
# combinations (Robo2 -> Robo2), replacing the original 2x2 (ECs / Neutro) grid.
$data = @(
    @("ECs","Robo2","Robo2","ECs",2,1,0.5215465,1.043093,0.3993854683363039,0.3071469942766215,2,1,0.5215465,1.043093,0.3993854683363039,0.3071469942766215,0.27201075166225,1.088043006649,0.1595087523182088,0.09433927609316294),
    @("ECs","Robo2","Robo2","FAPs",2,1,0.5215465,1.043093,0.3993854683363039,0.3071469942766215,1,0.3333333333333333,0.002342,0.007026,0.001793436954986034,0.002068861340060323,0.001221461903,0.007328771418,0.000716272658198732,0.0006354445421746313),
    @("ECs","Robo2","Robo2","Neutro",2,1,0.5215465,1.043093,0.3993854683363039,0.3071469942766215,3,1,0.7819839999999999,2.345952,0.59882109470871,0.6907841443833183,0.4078410182559999,2.447046109536,0.2391604433598963,0.2121722736412839),
    @("FAPs","Robo2","Robo2","ECs",1,0.3333333333333333,0.002342,0.007026,0.001793436954986034,0.002068861340060323,2,1,0.5215465,1.043093,0.3993854683363039,0.3071469942766215,0.001221461903,0.007328771418,0.000716272658198732,0.0006354445421746313),
    @("FAPs","Robo2","Robo2","FAPs",1,0.3333333333333333,0.002342,0.007026,0.001793436954986034,0.002068861340060323,1,0.3333333333333333,0.002342,0.007026,0.001793436954986034,0.002068861340060323,0.000005484964,0.000049364676,0.000003216416111509578,0.000004280187244396195),
    @("FAPs","Robo2","Robo2","Neutro",1,0.3333333333333333,0.002342,0.007026,0.001793436954986034,0.002068861340060323,3,1,0.7819839999999999,2.345952,0.59882109470871,0.6907841443833183,0.001831406528,0.016482658752,0.001073947880675793,0.001429136610641295),
    @("Neutro","Robo2","Robo2","ECs",3,1,0.7819839999999999,2.345952,0.59882109470871,0.6907841443833183,2,1,0.5215465,1.043093,0.3993854683363039,0.3071469942766215,0.4078410182559999,2.447046109536,0.2391604433598963,0.2121722736412839),
    @("Neutro","Robo2","Robo2","FAPs",3,1,0.7819839999999999,2.345952,0.59882109470871,0.6907841443833183,1,0.3333333333333333,0.002342,0.007026,0.001793436954986034,0.002068861340060323,0.001831406528,0.016482658752,0.001073947880675793,0.001429136610641295),
    @("Neutro","Robo2","Robo2","Neutro",3,1,0.7819839999999999,2.345952,0.59882109470871,0.6907841443833183,3,1,0.7819839999999999,2.345952,0.59882109470871,0.6907841443833183,0.6114989762559998,5.503490786303998,0.3585867034681379,0.4771827341313931)
)

$r = 2
foreach ($row in $data) {
    for ($c = 1; $c -le 20; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c - 1]
    }
    $r++
}
